$wb = $excel.ActiveWorkbook

# Overview sheet: row for 41711229-d0f4-47b9-aca9-e941cdd8c894.md (row 3)
# Status columns (zh-cn / de-de) change from "Handed back: in sync with en-US"
# to "Ready for handoff"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: same file's row (row 3) - Status updated, and a new handoff
# timestamp recorded for "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-20 03:28:24"

# de-de sheet: same file's row (row 3) - Status updated, and a new handoff
# timestamp recorded for "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-20 03:28:35"
